# Swap the words "ימין" (right) and "שמאל" (left) inside the "Details" paragraph,
# splitting them into their own runs (mirroring how Word splits a run when the
# user selects and retypes a word), and relocate the "_GoBack" bookmark to sit
# right after the new "שמאל" run (matching where the user's last edit landed),
# merging the run pair it used to separate back into one run.

$d = $word.ActiveDocument

# --- Step 1: swap the two words' text, without touching anything else -------
# Use a placeholder that cannot collide with any existing text so the two
# replacements don't clash with one another.
$placeholder = "@@SWAP_PLACEHOLDER@@"

$d.Content.Find.Execute("ימין", $true, $false, $false, $false, $false, $true, 1, $false, $placeholder, 2) | Out-Null
$d.Content.Find.Execute("שמאל", $true, $false, $false, $false, $false, $true, 1, $false, "ימין", 2) | Out-Null
$d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, "שמאל", 2) | Out-Null

# --- Step 2: force "שמאל" onto its own run ----------------------------------
$rngL = $d.Content
$rngL.Find.Execute("שמאל") | Out-Null
$rngL.Bold = $true
$rngL.Bold = $false

# --- Step 3: force "ימין" onto its own run -----------------------------------
$rngR = $d.Content
$rngR.Find.Execute("ימין") | Out-Null
$rngR.Bold = $true
$rngR.Bold = $false

# --- Step 4: move the _GoBack bookmark to just after the "שמאל" run ---------
$rngBm = $d.Content
$rngBm.Find.Execute("שמאל") | Out-Null
$newBmRange = $d.Range($rngBm.End, $rngBm.End)
$d.Bookmarks.Add("_GoBack", $newBmRange)

# --- Step 5: re-merge the runs that used to be split by the old bookmark ----
$d.Content.Find.Execute("סימולציה שנגמרה", $true, $false, $false, $false, $false, $true, 1, $false, "סימולציה שנגמרה", 2) | Out-Null
